# Login Details Updated changes
#
# The "login" worksheet had a duplicate/erroneous login entry in row 2
# (A2=9876543211) which needs to be removed. Deleting the entire row
# shifts all subsequent rows up by one (row 3's former value 987654321
# becomes the new row 2, etc.), which also changes the used
# range/dimension (A1:B16 -> A1:B15) and moves the mailto hyperlink that
# used to live on A8 up onto A7.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("login")

# Remember the hyperlink's target email address (the cell's own display
# text, e.g. "987654321@") before the row shift so it can be re-anchored
# at its new location afterwards.
$linkEmail = $ws.Range("A8").Text

# Delete row 2 entirely; remaining rows shift up automatically.
$ws.Rows.Item(2).Delete()

# The hyperlink object itself keeps pointing at the old A8 anchor, so
# drop it and re-create it on the cell it now belongs to (A7).
$ws.Range("A7").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A7"), "mailto:" + $linkEmail) | Out-Null

# Re-apply the built-in Hyperlink cell style so A7 keeps using the
# existing style rather than a newly duplicated one.
$ws.Range("A7").Style = "Hyperlink"

# Update the selected/active cell to reflect the post-edit state.
$ws.Range("E7").Select()
